$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Russell")
$srcDateCell = $wb.Worksheets.Item("Digital Ingest").Range("E4")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Examples"

$src.Range("A1:G2").Copy($ws.Range("A1:G2"))

$ws.Range("A3").Value = "Example"
$ws.Range("B3").Value = "Example"
$ws.Range("C3").Value = "Backlog"
$ws.Range("D3").Value = "E. X. Ample"
$srcDateCell.Copy($ws.Range("E3"))
$ws.Range("E3").Value = 45688
$ws.Range("F3").Value = "Info"

$ws.Range("A4").Value = "Example"
$ws.Range("B4").Value = "Example"
$ws.Range("C4").Value = "Backlog"
$srcDateCell.Copy($ws.Range("E4"))
$ws.Range("E4").Value = 45688
$ws.Range("F4").Value = "Info"

$ws.Range("C2").Validation.Add(3, 1, 1, '"Access/Mezzanine,Backlog,Outsourced Graphics,Medium Priority,Transfer,Working Files"')
$ws.Range("C2").Validation.IgnoreBlank = $false

"Active sheet: " + $wb.ActiveSheet.Name
"Russell tabselected check via Worksheets:"
$ws.Range("A1:G4").Select()
